$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-18 Sunday" "2026-01-19 Monday"

Replace-Text "26×38=" "61×89="
Replace-Text "48×99=" "36×53="
Replace-Text "32×98=" "32×16="
Replace-Text "54×97=" "80×72="
Replace-Text "90×89=" "79×23="
Replace-Text "49×39=" "76×79="
Replace-Text "64×71=" "62×71="
Replace-Text "56×71=" "31×31="
Replace-Text "48×27=" "71×92="
Replace-Text "17×78=" "81×62="
Replace-Text "64×95=" "97×49="
Replace-Text "80×14=" "82×16="
Replace-Text "90×84=" "15×94="
Replace-Text "38×31=" "25×21="
Replace-Text "50×61=" "34×54="
Replace-Text "75×90=" "33×62="
Replace-Text "67×69=" "19×47="
Replace-Text "12×27=" "13×82="
Replace-Text "99×50=" "78×54="
Replace-Text "86×71=" "94×75="
Replace-Text "70×16=" "78×51="
Replace-Text "26×47=" "93×90="
Replace-Text "96×84=" "66×61="
Replace-Text "22×21=" "13×76="
Replace-Text "61×63=" "13×25="
